# Fixed a bug in genDefaultScene
# The data rows (2-25) of the default scene table were being emitted in the
# wrong order; reorder them back to the correct sequence. Columns A-F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: new row number -> old row number (source of truth for that row's data)
$rowMap = @{
    2  = 8
    3  = 12
    4  = 15
    5  = 10
    6  = 13
    7  = 3
    8  = 9
    9  = 11
    10 = 14
    11 = 4
    12 = 5
    13 = 6
    14 = 7
    15 = 2
    16 = 17
    17 = 18
    18 = 21
    19 = 16
    20 = 19
    21 = 20
    22 = 23
    23 = 22
    24 = 24
    25 = 25
}

# Snapshot the original values for rows 2-25, columns A-F (1-6), before
# overwriting anything.
$original = @{}
for ($row = 2; $row -le 25; $row++) {
    $rowVals = @()
    for ($col = 1; $col -le 6; $col++) {
        $rowVals += , ($ws.Cells.Item($row, $col).Value2)
    }
    $original[$row] = $rowVals
}

# Write each target row using the values captured from its mapped source row.
foreach ($targetRow in $rowMap.Keys) {
    $sourceRow = $rowMap[$targetRow]
    $vals = $original[$sourceRow]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($targetRow, $col).Value2 = $vals[$col - 1]
    }
}
